$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns L, M, N: headers (reuse the same header style as the rest of row 1) ---
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)

# --- Rescale existing "particip" (E) and "taxa_sucesso" (F) columns from
#     fraction (0-1) to percentage (0-100) for rows 2..7 ---
$ws.Range("E2").Value = 97.37827715355806
$ws.Range("F2").Value = 62

$ws.Range("E3").Value = 2.621722846441948
$ws.Range("F3").Value = 68.57142857142857

$ws.Range("E4").Value = 95.50408719346049
$ws.Range("F4").Value = 94.00855920114122

$ws.Range("E5").Value = 4.49591280653951
$ws.Range("F5").Value = 98.48484848484848

$ws.Range("E6").Value = 98.9766081871345
$ws.Range("F6").Value = 22.15657311669129

$ws.Range("E7").Value = 1.023391812865497
$ws.Range("F7").Value = 28.57142857142857

# --- New data for columns L (apoio_medio), M (contribuicoes), N (media_contribuicoes) ---
$ws.Range("L2").Value = 91.56965423913746
$ws.Range("M2").Value = 255984
$ws.Range("N2").Value = 317.5980148883374

$ws.Range("L3").Value = 82.2979860710347
$ws.Range("M3").Value = 7569
$ws.Range("N3").Value = 315.375

$ws.Range("L4").Value = 91.21586592230445
$ws.Range("M4").Value = 185734
$ws.Range("N4").Value = 140.9210925644917

$ws.Range("L5").Value = 79.29010140385255
$ws.Range("M5").Value = 17912
$ws.Range("N5").Value = 275.5692307692308

$ws.Range("L6").Value = 19.6055125364595
$ws.Range("M6").Value = 2083
$ws.Range("N6").Value = 13.88666666666667

$ws.Range("L7").Value = 18.78940113071737
$ws.Range("M7").Value = 125
$ws.Range("N7").Value = 62.5
